$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.317.01"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "3.192.08"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'596.18"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").Value = "'154.06"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.190.00"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "'0.547"
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").Value = "'6.05"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'0.515"
$ws.Range("E12").Value = "  +4.53%  "
$ws.Range("D13").Value = "'0.0000266"
$ws.Range("E13").Value = "  +3.62%  "
$ws.Range("D14").Value = "'38.94"
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("D15").Value = "3.714.96"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "66.266.49"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("D17").Value = "'7.41"
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "3.183.88"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "'510.06"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("D22").Value = "'0.740"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("D23").Value = "'8.06"
$ws.Range("E23").Value = "  +5.87%  "
$ws.Range("D24").Value = "'15.03"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'85.18"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  +4.70%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("D30").Value = "'7.01"
$ws.Range("E30").Value = "  +14.22%  "
$ws.Range("D31").Value = "'2.89"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("D32").Value = "'28.18"
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'6.52"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "'54.88"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").Value = "'489.28"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").Value = "'0.0894"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'0.0422"
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("D40").Value = "'8.88"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.300"
$ws.Range("E41").Value = "  +7.00%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  +5.53%  "
$ws.Range("D43").Value = "'2.82"
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("D44").Value = "0.0₃0658"
$ws.Range("E44").Value = "  +15.35%  "
$ws.Range("D45").Value = "2.914.46"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'28.49"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'2.65"
$ws.Range("E50").Value = "  +11.83%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  +4.62%  "

Write-Host "Applied cryptos update"